$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric values in column D stay as text (avoid Excel auto-numeric conversion)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.942.65'
$ws.Range('E2').Value = '  +1.50%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.890.03'
$ws.Range('E3').Value = '  +1.26%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.03'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('E6').Value = '  -0.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4570'
$ws.Range('E7').Value = '  +0.68%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3904'
$ws.Range('E8').Value = '  +1.79%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07851'
$ws.Range('E9').Value = '  +0.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9884'
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.90'
$ws.Range('E11').Value = '  +1.86%  '
$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.037'
$ws.Range('E12').Value = '  +1.66%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.695'
$ws.Range('E13').Value = '  +1.04%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.772.02'
$ws.Range('E14').Value = '  -8.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06941'
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.09'
$ws.Range('E16').Value = '  +1.91%  '
$ws.Range('E17').Value = '  -0.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009994'
$ws.Range('E18').Value = '  +0.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.02'
$ws.Range('E19').Value = '  +2.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  -0.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '28.930.42'
$ws.Range('E21').Value = '  +1.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.296'
$ws.Range('E22').Value = '  +0.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.98'
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.078.16'
$ws.Range('E24').Value = '  -2.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.053'
$ws.Range('E25').Value = '  -1.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '156.12'
$ws.Range('E26').Value = '  +1.88%  '
$ws.Range('E27').Value = '  +1.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.943'
$ws.Range('E28').Value = '  +4.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.930'
$ws.Range('E29').Value = '  +2.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '117.72'
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09332'
$ws.Range('E31').Value = '  +0.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9104'
$ws.Range('E32').Value = '  +0.63%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.296'
$ws.Range('E33').Value = '  +0.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.328'
$ws.Range('E34').Value = '  +0.74%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.265'
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  +4.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05776'
$ws.Range('E37').Value = '  +2.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02067'
$ws.Range('E38').Value = '  +1.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.001'
$ws.Range('E39').Value = '  -0.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.644'
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5675'
$ws.Range('E41').Value = '  +1.97%  '
$ws.Range('E42').Value = '  +0.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.753'
$ws.Range('E43').Value = '  +1.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.306'
$ws.Range('E44').Value = '  +7.95%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '11.99'
$ws.Range('E45').Value = '  +3.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5366'
$ws.Range('E46').Value = '  +2.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.07045'
$ws.Range('E47').Value = '  -1.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.854'
$ws.Range('E48').Value = '  +2.74%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '112.87'
$ws.Range('E49').Value = '  +1.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.519'
$ws.Range('E50').Value = '  +3.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.080'
$ws.Range('E51').Value = '  -3.56%  '
